# Upload added data to the two tracker sheets: write a 0 into L1 on each
# sheet, update the remembered selection on each sheet, and leave the
# "assignmentTracker" sheet as the active/selected tab.

$wb  = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("emotionTracker")
$ws2 = $wb.Worksheets.Item("assignmentTracker")

# emotionTracker: value in L1, selection left on L1
$ws1.Range("L1").Value = 0
$ws1.Range("L1").Select()

# assignmentTracker: value in L1, selection left on F9
$ws2.Range("L1").Value = 0
$ws2.Range("F9").Select()

# assignmentTracker is the active sheet when the file is saved
$ws2.Activate()
